$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 82, shifting existing rows 82-85 down to 83-86
$ws.Rows.Item(82).Insert()

# Populate the new row 82 with the new record's data
$ws.Cells.Item(82, 1).Value = 1
$ws.Cells.Item(82, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(82, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(82, 4).Value = 44746
$ws.Cells.Item(82, 5).Value = 15
$ws.Cells.Item(82, 6).Value = 100112021
$ws.Cells.Item(82, 7).Value = "Ají"
$ws.Cells.Item(82, 8).Value = "Inferno"
$ws.Cells.Item(82, 9).Value = "Primera"
$ws.Cells.Item(82, 10).Value = 160
$ws.Cells.Item(82, 11).Value = 10000
$ws.Cells.Item(82, 12).Value = 11000
$ws.Cells.Item(82, 13).Value = 10500
$ws.Cells.Item(82, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(82, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(82, 16).Value = 700
$ws.Cells.Item(82, 17).Value = 15
$ws.Cells.Item(82, 18).Value = "Hortaliza"
